$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 28),
    @('Acquire BPO Southbank', 6),
    @('Adorn Cosmetics Clayton', 5),
    @('Al Haj Halal Meats Glenroy', 31),
    @('Al-Taqwa College Truganina', 9),
    @('Amiga Montessori Craigieburn', 28),
    @('Anglicare Home Princes Hwy Dandenong', 6),
    @('Best&Less Fountain Gate Narre Warren', 5),
    @('Budget Car and Truck Rentals Campbellfield', 7),
    @('CS Square Caroline Springs', 9),
    @('Campbellfield Ford Complex Vaccination Clinic Campbellfield', 12),
    @('Cannie Road Construction Site Cannie', 7),
    @('Caroline Springs Police Station', 11),
    @('Cedars Medical Clinic Coburg', 39),
    @('Chemist Warehouse Campbellfield DC', 6),
    @('Chemist Warehouse Fillo Drive Somerton', 7),
    @('City of Hobsons Bay Community', 5),
    @('City of Moreland Community', 5),
    @('City of Wyndham Community', 6),
    @('Classy Cabinets and Kitchens Craigieburn', 9),
    @('Coles Aurora Village Epping', 6),
    @('Coles Broadmeadows Central Shopping Centre', 8),
    @('Coles Campbellfield Plaza Campbellfield', 12),
    @('Coles Coburg North Village', 29),
    @('Coles Greenvale Shopping Centre', 7),
    @('Coles Pakenham Place Shopping', 11),
    @('Coles Roxburgh Village Roxburgh Park', 10),
    @('Community Kids Meadow Heights', 11),
    @('Construction Site Olea Apartment Caulfield North', 16),
    @('Costco Wholesale Epping', 28),
    @('Crusader Caravans Epping', 21),
    @('DRC Laverton Automotive Repairs Laverton North', 5),
    @('Direct Freight Express Campbellfield', 13),
    @('Don Watson Coldstore Derrimut', 5),
    @('Epworth Healthcare Epworth Richmond Emergency Department', 6),
    @('Fine Food Holdings Pty Ltd Dandenong South', 6),
    @('Fitzroy Community School Fitzroy North', 36),
    @('Fonterra Manufacturing Workplace Campbellfield', 7),
    @('General foods Campbellfield', 12),
    @('Glenroy West Primary School', 6),
    @('Goodstart Early Learning Altona', 7),
    @('Green Leaves Early Learning Centre Highlands Craigieburn', 8),
    @('Gumboots Early Learning Centre South Morang', 5),
    @('Hamilton Marino 236 Jasper Road McKinnon', 13),
    @('Health Care Providers Association South Melbourne', 9),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 6),
    @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 9),
    @('Ibis Kingsgate Hotel Melbourne', 6),
    @('Ilim College Glenroy Campus Hadfield', 16),
    @('Ilim Learning Sanctuary Glenroy', 11),
    @('Impact Designer Homes Epping', 5),
    @('Industrial Galvanizers Valmont Coatings Campbellfield', 22),
    @('KFC Fawkner', 7),
    @('Kasr Sweets Coolaroo', 5),
    @('Kids House Early Learning Cheltenham', 11),
    @('Level Crossing Removal Project Lilydale Construction Site John Street', 9),
    @('Lineage Logistics Laverton North', 8),
    @('Linfox Somerton National Distribution Centre Somerton', 9),
    @('McDonalds Thomastown II', 5),
    @('Melbourne Metropolitan Remand Centre Ravenhall', 11),
    @('Melbourne Truck Repairs Campbellfield', 7),
    @('Melbourne West Police Station Docklands', 7),
    @('Melbourne Youth Justice Centre Parkville', 5),
    @('Mercy Hospital for Women Heidelberg', 5),
    @('Mernda YMCA Early Learning Centre Mernda', 5),
    @('MyCentre Childcare Broadmeadows', 16),
    @('National Gallery of Victoria Melbourne', 9),
    @('Nido Early School Moonee Ponds', 15),
    @('North Geelong House Party', 7),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 48),
    @('Northern Health The Northern Hospital Epping', 8),
    @('OnQ Plumbing and Excavations Craigieburn', 18),
    @('Oporto Coolaroo', 12),
    @('Our Lady Help of Christian''s Primary School Brunswick East', 10),
    @('Paisley Park Early Learning Centre Bundoora', 7),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 13),
    @('Private Residence Northern Community Services Fawkner', 5),
    @('Ramsay Health Care Warringal Private Hospital Heidelberg', 9),
    @('Richmond Quarter 261-271 Bridge Road Construction Site', 11),
    @('Sacca''s Fruit World Broadmeadows Central Shopping Centre', 5),
    @('Salta Drive Construction Site Rangedale Drainage Altona North', 7),
    @('Sharpline Stainless Steel Coburg North', 5),
    @('St Margaret''s Primary School Maribyrnong', 11),
    @('St Vincents Hospital Emergency Department Melbourne', 5),
    @('Tek Foods Somerton', 14),
    @('The Homestead Child and Family Centre Roxburgh Park', 11),
    @('The Huntly-Goornong Rail Works', 5),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 17),
    @('ThorwestenCabinets Pakenham', 13),
    @('Total Window Concepts Hoppers Crossing', 6),
    @('Unilodge College Square Student Accommodation 570 Lygon Street', 14),
    @('Werribee Mercy Hospital Emergency Department', 10),
    @('Western Health Footscray Hospital Emergency Department', 6),
    @('Western Health Sunshine Hospital Emergency Department', 11),
    @('Woodlands Long Day Care and Kindergarten Roxburgh Park', 5),
    @('Woolworths Greenvale Lakes Roxburgh Park', 5),
    @('Yara Childcare Centre Truganina', 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
